$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1779.8
$ws.Range("J32").Value = 1779.8
$ws.Range("L32").Value = 1779.8
$ws.Range("N32").Value = -2431.8

$ws.Range("H43").Value = 1750
$ws.Range("I43").Value = 1833.3334
$ws.Range("K43").Value = 1833.3334
$ws.Range("M43").Value = -1764.3334

$ws.Range("H64").Value = 7844.222
$ws.Range("I64").Value = 3598
$ws.Range("J64").Value = 8375
$ws.Range("K64").Value = 3598
$ws.Range("L64").Value = 8375
$ws.Range("M64").Value = -3350
$ws.Range("N64").Value = -8871

$ws.Range("H67").Value = 7844.222
$ws.Range("I67").Value = 3598
$ws.Range("J67").Value = 8375
$ws.Range("K67").Value = 3598
$ws.Range("L67").Value = 8375
$ws.Range("M67").Value = -2740
$ws.Range("N67").Value = -10091

$ws.Range("H92").Value = 166.35715
$ws.Range("I92").Value = 94.09999999999999
$ws.Range("K92").Value = 94.09999999999999
$ws.Range("M92").Value = 1153.9

$ws.Range("H98").Value = 2842.9211
$ws.Range("I98").Value = 2628.8572
$ws.Range("J98").Value = 3442.3
$ws.Range("K98").Value = 2628.8572
$ws.Range("L98").Value = 3442.3
$ws.Range("M98").Value = -1130.8572
$ws.Range("N98").Value = -6438.3

$ws.Range("H122").Value = 2842.9211
$ws.Range("I122").Value = 2628.8572
$ws.Range("J122").Value = 3442.3
$ws.Range("K122").Value = 7886.571599999999
$ws.Range("L122").Value = 10326.9
$ws.Range("M122").Value = -5436.571599999999
$ws.Range("N122").Value = -15226.9

$ws.Range("H132").Value = 3264.9
$ws.Range("I132").Value = 3545.4443
$ws.Range("K132").Value = 10636.3329
$ws.Range("M132").Value = -8106.332900000001

$ws.Range("H137").Value = 3540.6562
$ws.Range("I137").Value = 4200.5
$ws.Range("J137").Value = 2089
$ws.Range("K137").Value = 12601.5
$ws.Range("L137").Value = 6267
$ws.Range("M137").Value = -10051.5
$ws.Range("N137").Value = -11367

$ws.Range("H138").Value = 4444.519
$ws.Range("J138").Value = 5792.0586
$ws.Range("L138").Value = 17376.1758
$ws.Range("N138").Value = -27656.1758

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2670.3635
$ws.Range("I45").Value = 2133.1667
$ws.Range("K45").Value = 2133.1667
$ws.Range("M45").Value = -1756.1667

$ws.Range("H61").Value = 4918.2856
$ws.Range("I61").Value = 3475
$ws.Range("K61").Value = 3475
$ws.Range("M61").Value = -3263

$ws.Range("H63").Value = 6072.1113
$ws.Range("I63").Value = 2662.25
$ws.Range("K63").Value = 2662.25
$ws.Range("M63").Value = -1976.25

$ws.Range("H66").Value = 6072.1113
$ws.Range("I66").Value = 2662.25
$ws.Range("K66").Value = 13311.25
$ws.Range("M66").Value = -9879.25

$ws.Range("H74").Value = 1940.3077
$ws.Range("I74").Value = 1829.5454
$ws.Range("K74").Value = 1829.5454
$ws.Range("M74").Value = -955.5454

$ws.Range("H77").Value = 1940.3077
$ws.Range("I77").Value = 1829.5454
$ws.Range("K77").Value = 9147.726999999999
$ws.Range("M77").Value = -4779.726999999999

$ws.Range("H110").Value = 136949.84
$ws.Range("I110").Value = 144701.25
$ws.Range("K110").Value = 144701.25
$ws.Range("M110").Value = -142656.25

$ws.Range("H132").Value = 4210.3687
$ws.Range("I132").Value = 4333.222
$ws.Range("J132").Value = 1999
$ws.Range("K132").Value = 12999.666
$ws.Range("L132").Value = 5997
$ws.Range("M132").Value = -10469.666
$ws.Range("N132").Value = -11057

$ws.Range("H136").Value = 4918.2856
$ws.Range("I136").Value = 3475
$ws.Range("K136").Value = 10425
$ws.Range("M136").Value = -7875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 44549.8
$ws.Range("J60").Value = 44549.8
$ws.Range("L60").Value = 44549.8
$ws.Range("N60").Value = -45747.8

$ws.Range("H105").Value = 63427.625
$ws.Range("I105").Value = 67556.13
$ws.Range("K105").Value = 67556.13
$ws.Range("M105").Value = -65809.13

$ws.Range("H134").Value = 117843.89
$ws.Range("I134").Value = 7574.375
$ws.Range("K134").Value = 22723.125
$ws.Range("M134").Value = -20188.125

$ws.Range("H137").Value = 59374.625
$ws.Range("J137").Value = 59374.625
$ws.Range("L137").Value = 59374.625
$ws.Range("N137").Value = -69574.625

$ws.Range("H139").Value = 79807.5
$ws.Range("J139").Value = 79807.5
$ws.Range("L139").Value = 79807.5
$ws.Range("N139").Value = -90087.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38544.242
$ws.Range("I31").Value = 2780.2
$ws.Range("J31").Value = 45995.082
$ws.Range("K31").Value = 2780.2
$ws.Range("L31").Value = 45995.082
$ws.Range("M31").Value = -2485.2
$ws.Range("N31").Value = -46585.082

$ws.Range("H34").Value = 38544.242
$ws.Range("I34").Value = 2780.2
$ws.Range("J34").Value = 45995.082
$ws.Range("K34").Value = 2780.2
$ws.Range("L34").Value = 45995.082
$ws.Range("M34").Value = -2578.2
$ws.Range("N34").Value = -46399.082

$ws.Range("H58").Value = 5032
$ws.Range("I58").Value = 4629.4683
$ws.Range("K58").Value = 4629.4683
$ws.Range("M58").Value = -4426.4683

$ws.Range("H132").Value = 1546.1428
$ws.Range("I132").Value = 1364.6
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4093.8
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1563.8
$ws.Range("N132").Value = -11060

$ws.Range("H134").Value = 296687.03
$ws.Range("I134").Value = 2647.2424
$ws.Range("K134").Value = 7941.7272
$ws.Range("M134").Value = -5406.7272

$ws.Range("H136").Value = 5032
$ws.Range("I136").Value = 4629.4683
$ws.Range("K136").Value = 13888.4049
$ws.Range("M136").Value = -11338.4049

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H92").Value = 1358.8889
$ws.Range("I92").Value = 682.5
$ws.Range("J92").Value = 1900
$ws.Range("K92").Value = 2047.5
$ws.Range("L92").Value = 5700
$ws.Range("M92").Value = -799.5
$ws.Range("N92").Value = -8196

$ws.Range("H131").Value = 3196.8667
$ws.Range("I131").Value = 1536.2858
$ws.Range("J131").Value = 3416.1887
$ws.Range("K131").Value = 4608.857400000001
$ws.Range("L131").Value = 10248.5661
$ws.Range("M131").Value = 431.1425999999992
$ws.Range("N131").Value = -20328.5661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 31253142
$ws.Range("I126").Value = 50002676
$ws.Range("K126").Value = 150008028
$ws.Range("M126").Value = -150005558

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 25017
$ws.Range("I30").Value = 30016
$ws.Range("K30").Value = 30016
$ws.Range("M30").Value = -29908

$ws.Range("H46").Value = 4405.2104
$ws.Range("I46").Value = 3823
$ws.Range("J46").Value = 5666.6665
$ws.Range("K46").Value = 3823
$ws.Range("L46").Value = 5666.6665
$ws.Range("M46").Value = -3635
$ws.Range("N46").Value = -6042.6665

$ws.Range("H132").Value = 6919.0713
$ws.Range("I132").Value = 4611
$ws.Range("K132").Value = 13833
$ws.Range("M132").Value = -11303

$ws.Range("H136").Value = 838798.5600000001
$ws.Range("I136").Value = 1116334
$ws.Range("J136").Value = 6192
$ws.Range("K136").Value = 3349002
$ws.Range("L136").Value = 18576
$ws.Range("M136").Value = -3346452
$ws.Range("N136").Value = -23676

$ws.Range("H137").Value = 55000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200

$ws.Range("H138").Value = 79995
$ws.Range("J138").Value = 79995
$ws.Range("L138").Value = 79995
$ws.Range("N138").Value = -90275

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9995
$ws.Range("J15").Value = 9995
$ws.Range("L15").Value = 9995
$ws.Range("N15").Value = -10571

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws.Range("H62").Value = 89332.914
$ws.Range("I62").Value = 204399
$ws.Range("K62").Value = 204399
$ws.Range("M62").Value = -203775

$ws.Range("H65").Value = 89332.914
$ws.Range("I65").Value = 204399
$ws.Range("K65").Value = 1021995
$ws.Range("M65").Value = -1018875

$ws.Range("H138").Value = 250066080
$ws.Range("J138").Value = 250066080
$ws.Range("L138").Value = 250066080
$ws.Range("N138").Value = -250076360

